# Add a new "Res.locals" list item right after the existing
# "Res.redirect(...)" list item, keeping the trailing _GoBack bookmark
# attached to the newly inserted (now last) paragraph - matching how
# Word itself behaves when you place the cursor at the very end of the
# last paragraph and press Enter to start a new list item.

$d = $word.ActiveDocument

# Locate the paragraph that currently ends the document ("Res.redirect...")
# via Find rather than a hard-coded offset, then expand to the whole
# paragraph so we know exactly where its paragraph mark sits.
$rng = $d.Content
$found = $rng.Find.Execute("Res.redirect(", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Res.redirect(' paragraph to anchor the new bullet on."
}
[void]$rng.Expand(4)  # wdParagraph - grow the hit to the full paragraph, incl. mark

# Position right before the paragraph mark (this is also exactly where the
# _GoBack bookmark currently sits).
$insertPos = $rng.End - 1
$newText = "Res.locals : lưu trữ các dl ở trong 1 vòng đời request response"

# Insert the new bullet's text first - this pushes the (zero-length)
# bookmark after the inserted text, still inside the original paragraph.
$textRange = $d.Range($insertPos, $insertPos)
$textRange.InsertBefore($newText)

# Now split the paragraph right before the freshly inserted text, which
# turns it into its own new list paragraph (inheriting the same
# ListParagraph style / numbering from the paragraph it was split off of)
# while leaving the bookmark attached after the new run, inside that new
# paragraph - exactly matching the target diff.
$splitRange = $d.Range($insertPos, $insertPos)
$splitRange.InsertParagraphBefore()
